# Notice of Entry Template -> Notice of Entry (pro template pack w/ placeholders)
# See commit: chore(templates): generate pro template pack + remove date
# placeholders from PDFs (#125)

$d = $word.ActiveDocument

# --- Paragraph 1: heading -------------------------------------------------
# Was bold/24pt-ish run-level formatting ("Notice of Entry Template");
# becomes a real Heading1-styled paragraph with the shorter title text.
$p1 = $d.Paragraphs.Item(1)
$p1.Style = "Heading1"
$r1 = $p1.Range
[void]$r1.MoveEnd(1, -1)
$r1.Text = "Notice of Entry"

# --- Paragraph 2: date line -> tenant name placeholder --------------------
$r2 = $d.Paragraphs.Item(2).Range
[void]$r2.MoveEnd(1, -1)
$r2.Text = "Tenant Name: {{TENANT_NAME}}"

# --- Paragraph 3: disclaimer (+ line break) -> property address -----------
# The old paragraph had a trailing <w:br/> after the disclaimer text; the
# replacement text (and removal of the break) is achieved by overwriting
# the whole paragraph range (minus its paragraph mark).
$r3 = $d.Paragraphs.Item(3).Range
[void]$r3.MoveEnd(1, -1)
$r3.Text = "Property Address: {{PROPERTY_ADDRESS}}"

# --- Paragraph 4: "Tenant:" -> unit number ---------------------------------
$r4 = $d.Paragraphs.Item(4).Range
[void]$r4.MoveEnd(1, -1)
$r4.Text = "Unit: {{UNIT_NUMBER}}"

# --- Paragraph 5: "Property:" -> date of notice ----------------------------
$r5 = $d.Paragraphs.Item(5).Range
[void]$r5.MoveEnd(1, -1)
$r5.Text = "Date of Notice: {{NOTICE_DATE}}"

# --- Paragraph 6: "Date of Entry:" -> planned entry date/time --------------
$r6 = $d.Paragraphs.Item(6).Range
[void]$r6.MoveEnd(1, -1)
$r6.Text = "Planned Entry Date/Time: {{ENTRY_DATE_TIME}}"

# --- Paragraph 7: "Reason:" -> reason for entry ----------------------------
$r7 = $d.Paragraphs.Item(7).Range
[void]$r7.MoveEnd(1, -1)
$r7.Text = "Reason for Entry: {{REASON_FOR_ENTRY}}"

# --- Paragraph 8: "Contact:" -> landlord/manager name ----------------------
$r8 = $d.Paragraphs.Item(8).Range
[void]$r8.MoveEnd(1, -1)
$r8.Text = "Landlord/Manager: {{LANDLORD_NAME}}"

Write-Host "Notice of Entry template updated to pro placeholder pack."
